$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark these invoices as Processed (column E) by setting the value to 1
$ws.Range("E27").Value = 1
$ws.Range("E32").Value = 1
$ws.Range("E33").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("E38").Value = 1
$ws.Range("E41").Value = 1
